# Updated cryptos list with latest prices and 1h volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.927.92'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").Value = '2.036.65'
$ws.Range("E3").Value = '  -0.85%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''228.12'
$ws.Range("E5").Value = '  -0.58%  '

$ws.Range("D6").Value = '''0.610'
$ws.Range("E6").Value = '  -0.77%  '

$ws.Range("D7").Value = '''60.92'
$ws.Range("E7").Value = '  +3.68%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Value = '''0.380'
$ws.Range("E9").Value = '  -1.54%  '

$ws.Range("D10").Value = '''0.0818'
$ws.Range("E10").Value = '  +0.74%  '

$ws.Range("E11").Value = '  +0.35%  '

$ws.Range("D12").Value = '2.337.67'
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("D13").Value = '''14.55'
$ws.Range("E13").Value = '  -0.40%  '

$ws.Range("D14").Value = '''21.40'
$ws.Range("E14").Value = '  +2.09%  '

$ws.Range("D15").Value = '''0.764'
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").Value = '''5.16'
$ws.Range("E16").Value = '  -2.06%  '

$ws.Range("D17").Value = '2.059.41'
$ws.Range("E17").Value = '  +1.04%  '

$ws.Range("D18").Value = '37.883.05'
$ws.Range("E18").Value = '  -0.20%  '

$ws.Range("D19").Value = '''69.78'
$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("D20").Value = '''5.92'
$ws.Range("E20").Value = '  -6.57%  '

$ws.Range("E21").Value = '  -1.34%  '

$ws.Range("D22").Value = '''224.33'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").Value = '''2.25'
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.36'
$ws.Range("E26").Value = '  +0.87%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''167.44'
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("D28").Value = '''0.130'
$ws.Range("E28").Value = '  -1.82%  '

$ws.Range("D29").Value = '''18.89'
$ws.Range("E29").Value = '  -0.60%  '

$ws.Range("E30").Value = '  -3.22%  '

$ws.Range("E31").Value = '  +0.91%  '

$ws.Range("D32").Value = '''2.27'
$ws.Range("E32").Value = '  +10.28%  '

$ws.Range("D33").Value = '''4.42'
$ws.Range("E33").Value = '  -2.66%  '

$ws.Range("D34").Value = '''0.0607'
$ws.Range("E34").Value = '  +0.26%  '

$ws.Range("D35").Value = '''4.51'
$ws.Range("E35").Value = '  -1.54%  '

$ws.Range("D36").Value = '''6.38'
$ws.Range("E36").Value = '  +5.13%  '

$ws.Range("D37").Value = '''2.29'
$ws.Range("E37").Value = '  -1.11%  '

$ws.Range("D38").Value = '''3.34'
$ws.Range("E38").Value = '  +2.18%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.540.37'
$ws.Range("E40").Value = '  +0.58%  '

$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''17.63'
$ws.Range("E41").Value = '  +6.32%  '

$ws.Range("E42").Value = '  +0.81%  '

$ws.Range("D43").Value = '''96.54'
$ws.Range("E43").Value = '  -1.02%  '

$ws.Range("E44").Value = '  -3.01%  '

$ws.Range("D45").Value = '''0.0915'
$ws.Range("E45").Value = '  -0.64%  '

$ws.Range("E46").Value = '  -2.17%  '

$ws.Range("D47").Value = '''4.03'
$ws.Range("E47").Value = '  -0.21%  '

$ws.Range("E48").Value = '  -0.57%  '

$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").Value = '''7.09'
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("D51").Value = '2.225.03'
$ws.Range("E51").Value = '  -0.86%  '
